$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.013.03"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "2.247.64"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "271.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.21%  "
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.630"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0972"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.93%  "
$ws.Range("E12").Value = "  +19.52%  "
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").Value = "2.582.82"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.69%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.817"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.25%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.246.37"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "43.950.37"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000106"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.88%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.26%  "
$ws.Range("E27").Value = "  +12.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0918"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.114"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.124"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0353"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +25.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "13.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("E41").Value = "  +14.18%  "
$ws.Range("E42").Value = "  +4.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.444"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "2.467.15"
$ws.Range("E51").Value = "  +2.66%  "
